$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header rich-text strings (Volume/Number and date range) ---
$volCell = $ws.Range("A8")
$volCell.Characters(21, 2).Text = "35"

$dateCell = $ws.Range("C9")
$dateCell.Characters(27, 9).Text = "8/28/2023"
$dateCell.Characters(47, 9).Text = "9/3/2023"

# --- Update weekly crime statistics table (rows 14-30) ---

# Row 14
$ws.Range("I14").Value = 48
$ws.Range("J14").Value = 47
$ws.Range("K14").Value = 2.127659574468
$ws.Range("L14").Value = -25
$ws.Range("M14").Value = -49.473684210526
$ws.Range("N14").Value = -85.498489425981

# Row 15
$ws.Range("C15").Value = 4
$ws.Range("E15").Value = -33.333333333333
$ws.Range("F15").Value = 14
$ws.Range("G15").Value = 25
$ws.Range("H15").Value = -44
$ws.Range("I15").Value = 153
$ws.Range("J15").Value = 166
$ws.Range("K15").Value = -7.831325301204
$ws.Range("L15").Value = 2.684563758389
$ws.Range("M15").Value = 3.378378378378
$ws.Range("N15").Value = -63.221153846153

# Row 16
$ws.Range("C16").Value = 46
$ws.Range("D16").Value = 59
$ws.Range("E16").Value = -22.033898305084
$ws.Range("F16").Value = 200
$ws.Range("G16").Value = 211
$ws.Range("H16").Value = -5.213270142180
$ws.Range("I16").Value = 1633
$ws.Range("J16").Value = 1755
$ws.Range("K16").Value = -6.951566951566
$ws.Range("L16").Value = 24.561403508771
$ws.Range("M16").Value = -29.368512110726
$ws.Range("N16").Value = -85.18820861678

# Row 17
$ws.Range("C17").Value = 65
$ws.Range("D17").Value = 86
$ws.Range("E17").Value = -24.418604651162
$ws.Range("F17").Value = 275
$ws.Range("G17").Value = 337
$ws.Range("H17").Value = -18.397626112759
$ws.Range("I17").Value = 2850
$ws.Range("J17").Value = 2827
$ws.Range("K17").Value = 0.813583303855
$ws.Range("L17").Value = 22.107969151671
$ws.Range("M17").Value = 24.72647702407
$ws.Range("N17").Value = -51.637536059731

# Row 18
$ws.Range("C18").Value = 41
$ws.Range("D18").Value = 37
$ws.Range("E18").Value = 10.810810810810
$ws.Range("F18").Value = 153
$ws.Range("G18").Value = 182
$ws.Range("H18").Value = -15.934065934065
$ws.Range("I18").Value = 1390
$ws.Range("J18").Value = 1618
$ws.Range("K18").Value = -14.091470951792
$ws.Range("L18").Value = 8.255451713395
$ws.Range("M18").Value = -33.460986117759
$ws.Range("N18").Value = -82.864891518737

# Row 19
$ws.Range("C19").Value = 103
$ws.Range("D19").Value = 134
$ws.Range("E19").Value = -23.134328358209
$ws.Range("F19").Value = 425
$ws.Range("G19").Value = 497
$ws.Range("H19").Value = -14.486921529175
$ws.Range("I19").Value = 3897
$ws.Range("J19").Value = 3984
$ws.Range("K19").Value = -2.183734939759
$ws.Range("L19").Value = 30.947580645161
$ws.Range("M19").Value = 40.179856115107
$ws.Range("N19").Value = -14.782418543625

# Row 20
$ws.Range("C20").Value = 31
$ws.Range("D20").Value = 46
$ws.Range("E20").Value = -32.608695652173
$ws.Range("F20").Value = 153
$ws.Range("G20").Value = 160
$ws.Range("H20").Value = -4.375
$ws.Range("I20").Value = 1217
$ws.Range("J20").Value = 1216
$ws.Range("K20").Value = 0.082236842105
$ws.Range("L20").Value = 21.457085828343
$ws.Range("M20").Value = 27.970557308096
$ws.Range("N20").Value = -80.685605459450

# Row 21
$ws.Range("C21").Value = 292
$ws.Range("D21").Value = 369
$ws.Range("E21").Value = -20.867208672086
$ws.Range("F21").Value = 1227
$ws.Range("G21").Value = 1414
$ws.Range("H21").Value = -13.224893917963
$ws.Range("I21").Value = 11188
$ws.Range("J21").Value = 11613
$ws.Range("K21").Value = -3.659691724791
$ws.Range("L21").Value = 22.675438596491
$ws.Range("M21").Value = 4.953095684803
$ws.Range("N21").Value = -69.474229898229

# Row 22
$ws.Range("C22").Value = 5
$ws.Range("D22").Value = 7
$ws.Range("E22").Value = -28.571428571428
$ws.Range("F22").Value = 18
$ws.Range("G22").Value = 21
$ws.Range("H22").Value = -14.285714285714
$ws.Range("I22").Value = 193
$ws.Range("J22").Value = 237
$ws.Range("K22").Value = -18.565400843881
$ws.Range("L22").Value = 21.383647798742
$ws.Range("M22").Value = -31.802120141342

# Row 23
$ws.Range("C23").Value = 19
$ws.Range("D23").Value = 34
$ws.Range("E23").Value = -44.117647058823
$ws.Range("F23").Value = 106
$ws.Range("G23").Value = 128
$ws.Range("H23").Value = -17.1875
$ws.Range("I23").Value = 1064
$ws.Range("J23").Value = 1048
$ws.Range("K23").Value = 1.526717557251
$ws.Range("L23").Value = 10.833333333333
$ws.Range("M23").Value = 35.887611749680

# Row 24
$ws.Range("C24").Value = 229
$ws.Range("D24").Value = 257
$ws.Range("E24").Value = -10.894941634241
$ws.Range("F24").Value = 986
$ws.Range("G24").Value = 1101
$ws.Range("H24").Value = -10.445049954586
$ws.Range("I24").Value = 8437
$ws.Range("J24").Value = 8938
$ws.Range("K24").Value = -5.605280823450
$ws.Range("L24").Value = 26.057074555505
$ws.Range("M24").Value = 22.648640790812

# Row 25
$ws.Range("C25").Value = 144
$ws.Range("D25").Value = 119
$ws.Range("E25").Value = 21.008403361344
$ws.Range("F25").Value = 484
$ws.Range("G25").Value = 430
$ws.Range("H25").Value = 12.558139534883
$ws.Range("I25").Value = 4162
$ws.Range("J25").Value = 4064
$ws.Range("K25").Value = 2.411417322834
$ws.Range("L25").Value = 35.702641017280
$ws.Range("M25").Value = -23.139427516158

# Row 26
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = -50
$ws.Range("F26").Value = 23
$ws.Range("G26").Value = 32
$ws.Range("H26").Value = -28.125
$ws.Range("I26").Value = 237
$ws.Range("J26").Value = 254
$ws.Range("K26").Value = -6.692913385826
$ws.Range("L26").Value = -6.324110671936

# Row 27
$ws.Range("C27").Value = 14
$ws.Range("D27").Value = 12
$ws.Range("E27").Value = 16.666666666666
$ws.Range("F27").Value = 67
$ws.Range("G27").Value = 52
$ws.Range("H27").Value = 28.846153846153
$ws.Range("I27").Value = 430
$ws.Range("J27").Value = 426
$ws.Range("K27").Value = 0.938967136150
$ws.Range("L27").Value = -8.315565031982

# Row 28
$ws.Range("C28").Value = 3
$ws.Range("D28").Value = 12
$ws.Range("E28").Value = -75
$ws.Range("F28").Value = 16
$ws.Range("G28").Value = 31
$ws.Range("H28").Value = -48.387096774193
$ws.Range("I28").Value = 156
$ws.Range("J28").Value = 247
$ws.Range("K28").Value = -36.842105263157
$ws.Range("L28").Value = -46.938775510204
$ws.Range("M28").Value = -56.424581005586
$ws.Range("N28").Value = -88.392857142857

# Row 29
$ws.Range("C29").Value = 3
$ws.Range("D29").Value = 7
$ws.Range("E29").Value = -57.142857142857
$ws.Range("F29").Value = 15
$ws.Range("G29").Value = 26
$ws.Range("H29").Value = -42.307692307692
$ws.Range("I29").Value = 137
$ws.Range("J29").Value = 206
$ws.Range("K29").Value = -33.495145631068
$ws.Range("L29").Value = -41.201716738197
$ws.Range("M29").Value = -52.758620689655
$ws.Range("N29").Value = -88.668320926385

# Row 30
$ws.Range("D30").Value = 2
$ws.Range("F30").Value = 1
$ws.Range("G30").Value = 8
$ws.Range("H30").Value = -87.5
$ws.Range("J30").Value = 53
$ws.Range("K30").Value = -22.641509433962
$ws.Range("L30").Value = -2.380952380952
